# ---------------------------------------------------------------------------
# WR_89708709_WeekEnding_072725.xlsx update
#
# - Refresh the "Report Generated On" timestamp
# - Zero out the Thursday pricing lines (H16, H17) and the summary total
#   (C8), bump the "Total Line Items" counter (C9) to 4
# - Replace the old row-18 TOTAL row with a new Thursday detail line
#   (Point 06 / ANC-DHM-10-84-D1) and push the Thursday TOTAL down to row 19
# - Append a brand-new "Friday (07/25/2025)" section (rows 22-25): section
#   header, column header row, one detail line (Point 01 / PLA-HDIG) and a
#   TOTAL row
# - Keep the merged-cell map in sync: drop A18:G18, add A19:G19, A22:H22
#   and A25:G25
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- small helpers ----------------------------------------------------------

function Set-RedBandStyle($range, $size, $align) {
    # Thick red banner style used for section headers (style s="7" / s="14")
    $range.Font.Name = "Calibri"
    $range.Font.Bold = $true
    $range.Font.Size = $size
    $range.Font.Color = 16777215   # white (BGR 0xFFFFFF)
    $range.Interior.Pattern = 1
    $range.Interior.Color = 192    # red  (BGR for RGB C00000)
    $range.HorizontalAlignment = $align
    $range.VerticalAlignment = -4108   # xlVAlignCenter
}

function Set-DetailTextStyle($range) {
    # Plain detail-row cell (style s="9")
    $range.Font.Name = "Calibri"
    $range.Font.Bold = $false
    $range.Font.Size = 11
}

function Set-DetailNumStyle($range) {
    # Numeric / unit-count detail cell, right aligned (style s="10")
    $range.Font.Name = "Calibri"
    $range.Font.Bold = $false
    $range.Font.Size = 11
    $range.HorizontalAlignment = -4152   # xlRight
}

function Set-DetailPriceStyle($range) {
    # Currency detail cell (style s="11")
    $range.Font.Name = "Calibri"
    $range.Font.Bold = $false
    $range.Font.Size = 11
    $range.HorizontalAlignment = -4152   # xlRight
    $range.NumberFormat = """$""#,##0.00_-"
}

function Set-TotalLabelStyle($range) {
    # Red "TOTAL" label cell (style s="15")
    $range.Font.Name = "Calibri"
    $range.Font.Bold = $true
    $range.Font.Size = 11
    $range.Font.Color = 16777215
    $range.Interior.Pattern = 1
    $range.Interior.Color = 192
    $range.HorizontalAlignment = -4152   # xlRight
}

function Set-TotalPriceStyle($range) {
    # Red TOTAL currency cell (style s="16")
    $range.Font.Name = "Calibri"
    $range.Font.Bold = $true
    $range.Font.Size = 11
    $range.Font.Color = 16777215
    $range.Interior.Pattern = 1
    $range.Interior.Color = 192
    $range.NumberFormat = """$""#,##0.00_-"
}

function Set-HeaderCellStyle($range) {
    # Red column-header cell, centered + wrapped (style s="8")
    $range.Font.Name = "Calibri"
    $range.Font.Bold = $true
    $range.Font.Size = 11
    $range.Font.Color = 16777215
    $range.Interior.Pattern = 1
    $range.Interior.Color = 192
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4108     # xlVAlignCenter
    $range.WrapText = $true
}

# --- 1. report-generated timestamp ------------------------------------------

$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"

# --- 2. summary numbers -------------------------------------------------

$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 4

# --- 3. zero out Thursday detail pricing ------------------------------------

$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0

# --- 4. unmerge the old TOTAL row before rewriting row 18 -------------------

$ws.Range("A18:G18").UnMerge()

# --- 5. row 18 becomes a new Thursday detail line ---------------------------

$ws.Range("A18").Value = "Point 06"
$ws.Range("B18").Value = "ANC-DHM-10-84-D1"
$ws.Range("C18").Value = "Inst"
$ws.Range("D18").Value = "ANC,Dbl Hlx Mach,10in,84in,Db Eye 1in"
$ws.Range("E18").Value = "EA"
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = 0

Set-DetailTextStyle  $ws.Range("A18:E18")
Set-DetailNumStyle   $ws.Range("F18:G18")
Set-DetailPriceStyle $ws.Range("H18")

# --- 6. row 19 is the (relocated) Thursday TOTAL row ------------------------

$ws.Range("A19").Value = "TOTAL"
$ws.Range("H19").Value = 0

Set-TotalLabelStyle $ws.Range("A19")
Set-TotalPriceStyle $ws.Range("H19")
$ws.Range("A19:G19").Merge()

# --- 7. new Friday section header (row 22) ----------------------------------

$ws.Range("A22").Value = "Friday (07/25/2025)"
Set-RedBandStyle $ws.Range("A22") 14 -4131   # xlLeft
$ws.Range("A22:H22").Merge()

# --- 8. new Friday column headers (row 23) ----------------------------------

$ws.Range("A23").Value = "Point Number"
$ws.Range("B23").Value = "Billable Unit Code"
$ws.Range("C23").Value = "Work Type"
$ws.Range("D23").Value = "Unit Description"
$ws.Range("E23").Value = "Unit of Measure"
$ws.Range("F23").Value = "# Units"
$ws.Range("G23").Value = "N/A"
$ws.Range("H23").Value = "Pricing"

Set-HeaderCellStyle $ws.Range("A23:H23")

# --- 9. new Friday detail line (row 24) -------------------------------------

$ws.Range("A24").Value = "Point 01"
$ws.Range("B24").Value = "PLA-HDIG"
$ws.Range("C24").Value = "Inst"
$ws.Range("D24").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("E24").Value = "EA"
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = 0

Set-DetailTextStyle  $ws.Range("A24:E24")
Set-DetailNumStyle   $ws.Range("F24:G24")
Set-DetailPriceStyle $ws.Range("H24")

# --- 10. new Friday TOTAL row (row 25) --------------------------------------

$ws.Range("A25").Value = "TOTAL"
$ws.Range("H25").Value = 0

Set-TotalLabelStyle $ws.Range("A25")
Set-TotalPriceStyle $ws.Range("H25")
$ws.Range("A25:G25").Merge()
